# Apply the "rotation" update to rows 2-7 of the active sheet.
# Two independent 3-cycles of row data occur:
#   rows 2,3,4 : row2<-old row3, row3<-old row4, row4<-old row2
#   rows 5,6,7 : row5<-old row6, row6<-old row7, row7<-old row5
# Columns A,B,D,E,F,G,H,I,J,P,Q,R carry the rotated data; all other
# columns are identical across the six rows so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7)

$A = @{ 2 = 111545414; 3 = 111545323; 4 = 111545401; 5 = 111545328; 6 = 111543957; 7 = 111543968 }
$B = @{ 2 = 57494;     3 = 57487;     4 = 57487;     5 = 57494;     6 = 57494;     7 = 57487 }
$D = @{ 2 = "LC"; 3 = "NT"; 4 = "NT"; 5 = "LC"; 6 = "LC"; 7 = "NT" }
$E = @{ 2 = 205992; 3 = 205998; 4 = 205998; 5 = 205992; 6 = 205992; 7 = 205998 }
$F = @{ 2 = "Vattenfladdermus"; 3 = "Nordfladdermus"; 4 = "Nordfladdermus"; 5 = "Vattenfladdermus"; 6 = "Vattenfladdermus"; 7 = "Nordfladdermus" }
$G = @{ 2 = "Myotis daubentonii"; 3 = "Eptesicus nilssonii"; 4 = "Eptesicus nilssonii"; 5 = "Myotis daubentonii"; 6 = "Myotis daubentonii"; 7 = "Eptesicus nilssonii" }
$H = @{ 2 = "(Kuhl, 1817)"; 3 = "(A.Keyserling & Blasius, 1839)"; 4 = "(A.Keyserling & Blasius, 1839)"; 5 = "(Kuhl, 1817)"; 6 = "(Kuhl, 1817)"; 7 = "(A.Keyserling & Blasius, 1839)" }
$I = @{ 2 = "9"; 3 = "2"; 4 = "6"; 5 = "1"; 6 = "1"; 7 = "256" }
$J = @{ 2 = "registreringar"; 3 = ""; 4 = "registreringar"; 5 = ""; 6 = "registreringar"; 7 = "" }
$P = @{
    2 = "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
    3 = "Orsa Viborg, intill en grupp med hålträd, Dlr"
    4 = "Orsa Viborg, glänta i mitten av skogsparti, Dlr"
    5 = "Orsa Viborg, intill en grupp med hålträd, Dlr"
    6 = "Orsa Viborg, glänta i skogsparti, Dlr"
    7 = "Orsa Viborg, glänta i skogsparti, Dlr"
}
$Q = @{ 2 = 480487.2503558649; 3 = 480427.8053356989; 4 = 480487.2503558649; 5 = 480427.8053356989; 6 = 480406.6045043401; 7 = 480406.6045043401 }
$R = @{ 2 = 6772784.264016891; 3 = 6772811.198980245; 4 = 6772784.264016891; 5 = 6772811.198980245; 6 = 6772745.04339793; 7 = 6772745.04339793 }

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = $A[$r]
    $ws.Cells.Item($r, 2).Value = $B[$r]
    $ws.Cells.Item($r, 4).Value = $D[$r]
    $ws.Cells.Item($r, 5).Value = $E[$r]
    $ws.Cells.Item($r, 6).Value = $F[$r]
    $ws.Cells.Item($r, 7).Value = $G[$r]
    $ws.Cells.Item($r, 8).Value = $H[$r]
    $ws.Cells.Item($r, 9).Value = $I[$r]
    $ws.Cells.Item($r, 10).Value = $J[$r]
    $ws.Cells.Item($r, 16).Value = $P[$r]
    $ws.Cells.Item($r, 17).Value = $Q[$r]
    $ws.Cells.Item($r, 18).Value = $R[$r]
}
